# Updated symbol list on Sat Dec 31 11:56:27 UTC 2022 with GitHub Actions
#
# Applies refreshed price/volume data to the crypto tracker sheet:
#  - Row 11..20 coin/link/rank cycle up one slot (LiechtensteinCryptoassetsExchange
#    moves from rank #10/row 11 down to rank #19/row 20, the rest shift up).
#  - Price (column D) values refreshed across many rows.
#  - Volume(1h) (column E) text refreshed to match new coin/rank pairing,
#    and a couple of Best/Worst-in-24h flags toggled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param(
        [string]$Addr,
        [string]$Val
    )
    # Column D ("Price") holds numeric-looking text (e.g. "246.56") that must
    # stay text, matching the source data's inline-string cell type. A
    # leading apostrophe is the standard Excel way to force text storage for
    # a numeric-looking entry without altering the cell's number format.
    $ws.Range($Addr).Value = "'" + $Val
}

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    # Columns B/C/E hold non-numeric-looking text already, so a plain value
    # assignment keeps the cell text without needing a format change.
    $ws.Range($Addr).Value = $Val
}

# ---- Simple price (column D) refreshes ----
Set-PriceText "D2" "246.56"
Set-PriceText "D4" "5.098"
Set-PriceText "D5" "0.05598"
Set-PriceText "D6" "6.502"
Set-PriceText "D8" "0.8109"
Set-PriceText "D9" "0.8406"
Set-PriceText "D10" "0.1347"

# ---- Rows 11-20: coin/link/rank cycle up one row, prices refreshed ----
Set-TextValue "B11" "BitrueCoin"
Set-TextValue "C11" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-PriceText "D11" "0.02848"
Set-TextValue "E11" "10BitrueCoinBTR"

Set-TextValue "B12" "BitMartToken"
Set-TextValue "C12" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-PriceText "D12" "0.09404"
Set-TextValue "E12" "11BitMartTokenBMX"

Set-TextValue "B13" "BitForexToken"
Set-TextValue "C13" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-PriceText "D13" "0.001512"
Set-TextValue "E13" "12BitForexTokenBF"

Set-TextValue "B14" "One"
Set-TextValue "C14" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-PriceText "D14" "0.0006010"
Set-TextValue "E14" "13OneONE"

Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-PriceText "D15" "0.006120"
Set-TextValue "E15" "14TigerCashTCH"

Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-PriceText "D16" "3.556"
Set-TextValue "E16" "15LEOLEO"

Set-TextValue "B17" "BTSEToken"
Set-TextValue "C17" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-PriceText "D17" "2.118"
Set-TextValue "E17" "16BTSETokenBTSE"

Set-TextValue "B18" "BitpandaEcosystemToken"
Set-TextValue "C18" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-PriceText "D18" "0.3183"
Set-TextValue "E18" "17BitpandaEcosystemTokenBEST"

Set-TextValue "B19" "MandalaExchangeToken"
Set-TextValue "C19" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-PriceText "D19" "0.07007"
Set-TextValue "E19" "18MandalaExchangeTokenMDX"

Set-TextValue "B20" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C20" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-PriceText "D20" "0.03229"
Set-TextValue "E20" "19LiechtensteinCryptoassetsExchangeLCX"

# ---- More simple price (column D) refreshes ----
Set-PriceText "D22" "3.742"
Set-PriceText "D23" "0.04686"
Set-PriceText "D25" "0.001246"
Set-PriceText "D26" "0.004606"
Set-PriceText "D27" "0.00009599"

Set-PriceText "D41" "0.006110"
Set-TextValue "E41" "40KickTokenKICKBestin24h"

Set-PriceText "D42" "0.1053"
Set-PriceText "D43" "0.002500"

Set-PriceText "D44" "0.008711"
Set-TextValue "E44" "43LocalTradersLCT"

Set-PriceText "D45" "0.00005291"

Set-TextValue "E47" "46CoinbaseStockTokenCOINWorstin24h"

Set-PriceText "D48" "0.002053"
